# C5-PowerPoint.pptx edit
#
# 1) Slide 6 contains a 4-column table (graphicFrame) whose table style
#    was changed from {D1C8E043-8ED5-4074-B2CB-3B7D307C7D48} to
#    {82950C7E-B4E5-444B-A495-40B77E045FFB}.
#
# 2) The deck's "Office Theme" / "Integral" theme color palettes were
#    swapped between the two theme parts backing the deck. We reproduce
#    the reachable part of that swap by re-pointing the live theme's
#    color slots (Office Theme's palette) via the per-slide
#    ThemeColorScheme object exposed on the object model.

function ToOle([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1) table style on slide 6 -------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{82950C7E-B4E5-444B-A495-40B77E045FFB}")
    }
}

# --- 2) theme color palette -----------------------------------------------------
# Target ("Office Theme") srgb values, in clrScheme slot order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
#   7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    (ToOle 0x00 0x00 0x00),  # dk1      000000
    (ToOle 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (ToOle 0x44 0x54 0x6A),  # dk2      44546A
    (ToOle 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (ToOle 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (ToOle 0xED 0x7D 0x31),  # accent2  ED7D31
    (ToOle 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (ToOle 0xFF 0xC0 0x00),  # accent4  FFC000
    (ToOle 0x44 0x72 0xC4),  # accent5  4472C4
    (ToOle 0x70 0xAD 0x47),  # accent6  70AD47
    (ToOle 0x05 0x63 0xC1),  # hlink    0563C1
    (ToOle 0x95 0x4F 0x72)   # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
